# Insert a new weekly price record as row 41, pushing the existing
# rows 41-165 down to 42-166 (row 166 ends up holding what used to be
# row 165, unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(41).Insert()

$ws.Range("A41").Value = 3
$ws.Range("B41").Value = "Femacal de La Calera"
$ws.Range("C41").Value = "Coquimbo"
$ws.Range("D41").Value = 44648
$ws.Range("E41").Value = 5
$ws.Range("F41").Value = 100112052
$ws.Range("G41").Value = "Albahaca"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 90
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 4800
$ws.Range("M41").Value = 4633
$ws.Range("N41").Value = "$/docena de matas"
$ws.Range("O41").Value = "Provincia de Quillota"
$ws.Range("P41").Value = 772
$ws.Range("Q41").Value = 6
$ws.Range("R41").Value = "Hortaliza"
